$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.079.22'
$ws.Range("E2").Value = '  +4.58%  '

$ws.Range("D3").Value = '3.254.50'
$ws.Range("E3").Value = '  +2.84%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("E5").Value = '  +2.85%  '

$ws.Range("D6").Value = '''177.42'
$ws.Range("E6").Value = '  +4.48%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '''0.602'
$ws.Range("E8").Value = '  -0.93%  '

$ws.Range("D9").Value = '3.252.69'
$ws.Range("E9").Value = '  +2.93%  '

$ws.Range("E10").Value = '  +4.20%  '

$ws.Range("E11").Value = '  +2.03%  '

$ws.Range("E12").Value = '  +4.49%  '

$ws.Range("D13").Value = '3.814.89'
$ws.Range("E13").Value = '  +2.77%  '

$ws.Range("E14").Value = '  +0.52%  '

$ws.Range("D15").Value = '''28.10'
$ws.Range("E15").Value = '  +2.95%  '

$ws.Range("D16").Value = '67.069.63'
$ws.Range("E16").Value = '  +4.66%  '

$ws.Range("E17").Value = '  +3.07%  '

$ws.Range("D18").Value = '3.257.67'
$ws.Range("E18").Value = '  +3.07%  '

$ws.Range("E19").Value = '  +2.29%  '

$ws.Range("D20").Value = '''13.45'
$ws.Range("E20").Value = '  +3.35%  '

$ws.Range("D21").Value = '''372.72'
$ws.Range("E21").Value = '  +5.35%  '

$ws.Range("D22").Value = '''7.63'
$ws.Range("E22").Value = '  +5.81%  '

$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").Value = '''71.04'
$ws.Range("E24").Value = '  +2.75%  '

$ws.Range("D25").Value = '''0.512'
$ws.Range("E25").Value = '  +1.69%  '

$ws.Range("D26").Value = '3.397.02'
$ws.Range("E26").Value = '  +2.75%  '

$ws.Range("D27").Value = '''0.0000119'
$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("D28").Value = '''9.82'
$ws.Range("E28").Value = '  +2.36%  '

$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  +4.39%  '

$ws.Range("D32").Value = '''5.61'
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").Value = '''22.59'
$ws.Range("E33").Value = '  +2.33%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").Value = '''1.27'
$ws.Range("E35").Value = '  +5.98%  '

$ws.Range("E36").Value = '  +2.80%  '

$ws.Range("D37").Value = '''167.33'
$ws.Range("E37").Value = '  +8.41%  '

$ws.Range("E38").Value = '  +4.82%  '

$ws.Range("D39").Value = '''0.854'
$ws.Range("E39").Value = '  +5.64%  '

$ws.Range("D40").Value = '''1.88'
$ws.Range("E40").Value = '  +10.39%  '

$ws.Range("D41").Value = '''27.18'
$ws.Range("E41").Value = '  +5.26%  '

$ws.Range("D42").Value = '''2.59'
$ws.Range("E42").Value = '  +1.85%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''6.55'
$ws.Range("E43").Value = '  +9.24%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.758.38'
$ws.Range("E44").Value = '  +5.97%  '

$ws.Range("D45").Value = '''353.80'
$ws.Range("E45").Value = '  +9.41%  '

$ws.Range("E46").Value = '  +4.85%  '

$ws.Range("E47").Value = '  +5.91%  '

$ws.Range("D48").Value = '''40.50'
$ws.Range("E48").Value = '  +2.79%  '

$ws.Range("E49").Value = '  +2.64%  '

$ws.Range("D50").Value = '''0.0279'
$ws.Range("E50").Value = '  +3.26%  '

$ws.Range("E51").Value = '  +0.63%  '
